$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-29 Tuesday" "2023-08-30 Wednesday"

Replace-Text "98×59=5782" "86×84=7224"
Replace-Text "52×13=676" "41×23=943"
Replace-Text "51×63=3213" "42×61=2562"
Replace-Text "43×25=1075" "16×42=672"
Replace-Text "82×32=2624" "77×76=5852"
Replace-Text "30×58=1740" "23×64=1472"
Replace-Text "39×52=2028" "67×85=5695"
Replace-Text "95×18=1710" "64×58=3712"
Replace-Text "31×30=930" "47×43=2021"
Replace-Text "14×16=224" "69×99=6831"
Replace-Text "69×28=1932" "48×87=4176"
Replace-Text "94×19=1786" "40×92=3680"
Replace-Text "64×16=1024" "16×15=240"
Replace-Text "40×88=3520" "92×96=8832"
Replace-Text "23×13=299" "48×53=2544"
Replace-Text "97×22=2134" "59×56=3304"
Replace-Text "40×69=2760" "25×87=2175"
Replace-Text "57×25=1425" "25×18=450"
Replace-Text "70×59=4130" "73×46=3358"
Replace-Text "33×73=2409" "27×56=1512"
Replace-Text "51×84=4284" "35×75=2625"
Replace-Text "88×26=2288" "29×90=2610"
Replace-Text "56×84=4704" "58×86=4988"
Replace-Text "47×96=4512" "64×35=2240"
Replace-Text "89×30=2670" "68×12=816"
